$d = $word.ActiveDocument

# The "Prato, data emissione" paragraph was missing the {{data_odierna}}
# replacement placeholder - it had literal text "data emissione" instead
# of the template tag. Restore it as "Prato, {{data_odierna}}".
$d.Content.Find.Execute("Prato, data emissione", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Prato, {{data_odierna}}", 2)
